# Manual dislocation uploading 2021/08/21 13:00
# Shift the August-2021 dislocation data forward by one month (31 days) to
# September 2021, update the a few CarAmount figures that changed with the
# new upload, and clear out the two now-stale trailing rows (32 and 33)
# that used to hold next-month / duplicate data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift ShippingDate (column A) from Aug-2021 to Sep-2021, row by row,
#     and update CarAmount (column B) where the reloaded figures differ ---
$ws.Range("A2").Value  = 44440
$ws.Range("A3").Value  = 44441
$ws.Range("A4").Value  = 44442
$ws.Range("A5").Value  = 44443
$ws.Range("A6").Value  = 44444

$ws.Range("A7").Value  = 44445
$ws.Range("B7").Value  = 17

$ws.Range("A8").Value  = 44446
$ws.Range("A9").Value  = 44447

$ws.Range("A10").Value = 44448
$ws.Range("B10").Value = 0

$ws.Range("A11").Value = 44449
$ws.Range("B11").Value = 9

$ws.Range("A12").Value = 44450
$ws.Range("A13").Value = 44451
$ws.Range("A14").Value = 44452
$ws.Range("A15").Value = 44453
$ws.Range("A16").Value = 44454

$ws.Range("A17").Value = 44455
$ws.Range("B17").Value = 20

$ws.Range("A18").Value = 44456
$ws.Range("A19").Value = 44457
$ws.Range("A20").Value = 44458

$ws.Range("A21").Value = 44459
$ws.Range("B21").Value = 0

$ws.Range("A22").Value = 44460
$ws.Range("A23").Value = 44461

$ws.Range("A24").Value = 44462
$ws.Range("B24").Value = 17

$ws.Range("A25").Value = 44463
$ws.Range("A26").Value = 44464
$ws.Range("A27").Value = 44465

$ws.Range("A28").Value = 44466
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = 44467

$ws.Range("A30").Value = 44468
$ws.Range("B30").Value = 31

$ws.Range("A31").Value = 44469

# --- Rows 32 and 33 no longer carry data after the shift; clear them out
#     (keeps the existing per-cell formatting / styles, just like Delete
#     key on a selection in Excel) ---
$ws.Range("A32:E32").ClearContents()
$ws.Range("A33:E33").ClearContents()

# --- Update the sheet selection / view to match where the upload left
#     the cursor ---
$ws.Range("B2:B31").Select()
